# Fruta / hortaliza, semanal
# Insert a new weekly record as row 80 (pushing the existing rows 80-89
# down to 81-90) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 80:89 down one position, leaving a blank row 80 to fill in.
$ws.Rows("80:80").Insert()

$ws.Range("A80").Value2 = 9
$ws.Range("B80").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C80").Value2 = "Metropolitana"
$ws.Range("D80").Value2 = 44826
$ws.Range("E80").Value2 = 13
$ws.Range("F80").Value2 = "Fruta"
$ws.Range("G80").Value2 = 100102
$ws.Range("H80").Value2 = "Cítricos"
$ws.Range("I80").Value2 = 100102006
$ws.Range("J80").Value2 = "Pomelo"
$ws.Range("K80").Value2 = "Start Ruby"
$ws.Range("L80").Value2 = "Primera"
$ws.Range("M80").Value2 = 300
$ws.Range("N80").Value2 = 12000
$ws.Range("O80").Value2 = 12000
$ws.Range("P80").Value2 = 12000
$ws.Range("Q80").Value2 = "`$/caja 14 kilos"
$ws.Range("R80").Value2 = "Región Metropolitana"
$ws.Range("S80").Value2 = 857
$ws.Range("T80").Value2 = 14
